$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 4.426473626474878
$ws.Cells.Item(2, 5).Value = 5.000070571899414
$ws.Cells.Item(2, 6).Value = 5.0683119285427
$ws.Cells.Item(2, 7).Value = 4.402496874177151
$ws.Cells.Item(2, 8).Value = 262773412
$ws.Cells.Item(2, 9).Value = "TTEK"

$ws.Cells.Item(3, 4).Value = 4.776673203384721
$ws.Cells.Item(3, 5).Value = 4.928372383117676
$ws.Cells.Item(3, 6).Value = 4.989422026115578
$ws.Cells.Item(3, 7).Value = 4.462174828633395
$ws.Cells.Item(3, 8).Value = 262773412
$ws.Cells.Item(3, 9).Value = "TTEK"

$ws.Cells.Item(4, 4).Value = 4.534691456743436
$ws.Cells.Item(4, 5).Value = 4.991129398345947
$ws.Cells.Item(4, 6).Value = 5.111733000422605
$ws.Cells.Item(4, 7).Value = 4.415943600955563
$ws.Cells.Item(4, 8).Value = 262773412
$ws.Cells.Item(4, 9).Value = "TTEK"

$ws.Cells.Item(5, 4).Value = 4.799152438320951
$ws.Cells.Item(5, 5).Value = 4.929412364959717
$ws.Cells.Item(5, 6).Value = 5.061533458257486
$ws.Cells.Item(5, 7).Value = 4.252060391507659
$ws.Cells.Item(5, 8).Value = 262773412
$ws.Cells.Item(5, 9).Value = "TTEK"

$ws.Cells.Item(6, 4).Value = 5.497273789298037
$ws.Cells.Item(6, 5).Value = 5.487940788269043
$ws.Cells.Item(6, 6).Value = 5.924736006471511
$ws.Cells.Item(6, 7).Value = 5.271409467843365
$ws.Cells.Item(6, 8).Value = 262773412
$ws.Cells.Item(6, 9).Value = "TTEK"

$ws.Cells.Item(7, 4).Value = 5.742373769169701
$ws.Cells.Item(7, 5).Value = 6.165515899658203
$ws.Cells.Item(7, 6).Value = 6.43325642634913
$ws.Cells.Item(7, 7).Value = 5.633780087110474
$ws.Cells.Item(7, 8).Value = 262773412
$ws.Cells.Item(7, 9).Value = "TTEK"

$ws.Cells.Item(8, 4).Value = 6.6307820765162
$ws.Cells.Item(8, 5).Value = 7.218390941619873
$ws.Cells.Item(8, 6).Value = 7.36857838513435
$ws.Cells.Item(8, 7).Value = 6.529405194068817
$ws.Cells.Item(8, 8).Value = 262773412
$ws.Cells.Item(8, 9).Value = "TTEK"

$ws.Cells.Item(9, 4).Value = 8.164549889794488
$ws.Cells.Item(9, 5).Value = 8.220986366271973
$ws.Cells.Item(9, 6).Value = 8.437328584879983
$ws.Cells.Item(9, 7).Value = 7.713053592766514
$ws.Cells.Item(9, 8).Value = 262773412
$ws.Cells.Item(9, 9).Value = "TTEK"

$ws.Cells.Item(10, 4).Value = 7.710943827444272
$ws.Cells.Item(10, 5).Value = 8.285964965820312
$ws.Cells.Item(10, 6).Value = 8.427363326391664
$ws.Cells.Item(10, 7).Value = 7.522411930858428
$ws.Cells.Item(10, 8).Value = 262773412
$ws.Cells.Item(10, 9).Value = "TTEK"

$ws.Cells.Item(11, 4).Value = 8.691390080031109
$ws.Cells.Item(11, 5).Value = 8.965357780456543
$ws.Cells.Item(11, 6).Value = 9.031488536387997
$ws.Cells.Item(11, 7).Value = 8.49299961414016
$ws.Cells.Item(11, 8).Value = 262773412
$ws.Cells.Item(11, 9).Value = "TTEK"

$ws.Cells.Item(12, 4).Value = 8.816452425038126
$ws.Cells.Item(12, 5).Value = 9.327825546264648
$ws.Cells.Item(12, 6).Value = 9.46987328428188
$ws.Cells.Item(12, 7).Value = 8.72175333094831
$ws.Cells.Item(12, 8).Value = 262773412
$ws.Cells.Item(12, 9).Value = "TTEK"

$ws.Cells.Item(13, 4).Value = 9.166139723934114
$ws.Cells.Item(13, 5).Value = 9.431824684143066
$ws.Cells.Item(13, 6).Value = 9.659555554954901
$ws.Cells.Item(13, 7).Value = 9.014319746671854
$ws.Cells.Item(13, 8).Value = 262773412
$ws.Cells.Item(13, 9).Value = "TTEK"

$ws.Cells.Item(14, 4).Value = 9.347551857379869
$ws.Cells.Item(14, 5).Value = 9.204914093017578
$ws.Cells.Item(14, 6).Value = 10.02270569530771
$ws.Cells.Item(14, 7).Value = 9.111723349380933
$ws.Cells.Item(14, 8).Value = 262773412
$ws.Cells.Item(14, 9).Value = "TTEK"

$ws.Cells.Item(15, 4).Value = 11.07503381067518
$ws.Cells.Item(15, 5).Value = 11.58970832824707
$ws.Cells.Item(15, 6).Value = 11.86610758936998
$ws.Cells.Item(15, 7).Value = 11.0559722628122
$ws.Cells.Item(15, 8).Value = 262773412
$ws.Cells.Item(15, 9).Value = "TTEK"

$ws.Cells.Item(16, 4).Value = 13.17527122257567
$ws.Cells.Item(16, 5).Value = 12.61007118225098
$ws.Cells.Item(16, 6).Value = 13.37385519278459
$ws.Cells.Item(16, 7).Value = 12.15370973290519
$ws.Cells.Item(16, 8).Value = 262773412
$ws.Cells.Item(16, 9).Value = "TTEK"

$ws.Cells.Item(17, 4).Value = 9.724654976617645
$ws.Cells.Item(17, 5).Value = 10.55879783630371
$ws.Cells.Item(17, 6).Value = 10.62001947287349
$ws.Cells.Item(17, 7).Value = 9.282712285426626
$ws.Cells.Item(17, 8).Value = 262773412
$ws.Cells.Item(17, 9).Value = "TTEK"

$ws.Cells.Item(18, 4).Value = 11.49564799398088
$ws.Cells.Item(18, 5).Value = 12.40824508666992
$ws.Cells.Item(18, 6).Value = 12.43316896261018
$ws.Cells.Item(18, 7).Value = 11.27900234703924
$ws.Cells.Item(18, 8).Value = 262773412
$ws.Cells.Item(18, 9).Value = "TTEK"

$ws.Cells.Item(19, 4).Value = 15.31512588161407
$ws.Cells.Item(19, 5).Value = 15.2190465927124
$ws.Cells.Item(19, 6).Value = 16.63334193529208
$ws.Cells.Item(19, 7).Value = 15.03841675989472
$ws.Cells.Item(19, 8).Value = 262773412
$ws.Cells.Item(19, 9).Value = "TTEK"

$ws.Cells.Item(20, 4).Value = 16.68043795644946
$ws.Cells.Item(20, 5).Value = 16.84023475646973
$ws.Cells.Item(20, 6).Value = 17.39278518150533
$ws.Cells.Item(20, 7).Value = 15.63117381569685
$ws.Cells.Item(20, 8).Value = 262773412
$ws.Cells.Item(20, 9).Value = "TTEK"

$ws.Cells.Item(21, 4).Value = 16.70263910259202
$ws.Cells.Item(21, 5).Value = 16.50786209106445
$ws.Cells.Item(21, 6).Value = 17.47789045131807
$ws.Cells.Item(21, 7).Value = 16.20123165964487
$ws.Cells.Item(21, 8).Value = 262773412
$ws.Cells.Item(21, 9).Value = "TTEK"

$ws.Cells.Item(22, 4).Value = 13.1208503874493
$ws.Cells.Item(22, 5).Value = 14.54052066802978
$ws.Cells.Item(22, 6).Value = 16.093466651625
$ws.Cells.Item(22, 7).Value = 12.49117279435246
$ws.Cells.Item(22, 8).Value = 262773412
$ws.Cells.Item(22, 9).Value = "TTEK"

$ws.Cells.Item(23, 4).Value = 15.39023354701433
$ws.Cells.Item(23, 5).Value = 17.16372108459473
$ws.Cells.Item(23, 6).Value = 17.95172378726301
$ws.Cells.Item(23, 7).Value = 14.51317014818697
$ws.Cells.Item(23, 8).Value = 262773412
$ws.Cells.Item(23, 9).Value = "TTEK"

$ws.Cells.Item(24, 4).Value = 18.57752448188415
$ws.Cells.Item(24, 5).Value = 19.57254028320312
$ws.Cells.Item(24, 6).Value = 21.83217853310516
$ws.Cells.Item(24, 7).Value = 18.48442285262029
$ws.Cells.Item(24, 8).Value = 262773412
$ws.Cells.Item(24, 9).Value = "TTEK"

$ws.Cells.Item(25, 4).Value = 22.64934460424917
$ws.Cells.Item(25, 5).Value = 23.61273193359375
$ws.Cells.Item(25, 6).Value = 27.8430965402503
$ws.Cells.Item(25, 7).Value = 21.97730204612043
$ws.Cells.Item(25, 8).Value = 262773412
$ws.Cells.Item(25, 9).Value = "TTEK"

$ws.Cells.Item(26, 4).Value = 26.74018856046581
$ws.Cells.Item(26, 5).Value = 24.82072830200196
$ws.Cells.Item(26, 6).Value = 28.15401571739771
$ws.Cells.Item(26, 7).Value = 24.42400207555731
$ws.Cells.Item(26, 8).Value = 262773412
$ws.Cells.Item(26, 9).Value = "TTEK"

$ws.Cells.Item(27, 4).Value = 23.95422959132552
$ws.Cells.Item(27, 5).Value = 26.00934219360352
$ws.Cells.Item(27, 6).Value = 26.08141663227421
$ws.Cells.Item(27, 7).Value = 23.45749613428396
$ws.Cells.Item(27, 8).Value = 262773412
$ws.Cells.Item(27, 9).Value = "TTEK"

$ws.Cells.Item(28, 4).Value = 29.38432136346508
$ws.Cells.Item(28, 5).Value = 34.26935195922852
$ws.Cells.Item(28, 6).Value = 34.37079838255507
$ws.Cells.Item(28, 7).Value = 28.67224472028209
$ws.Cells.Item(28, 8).Value = 262773412
$ws.Cells.Item(28, 9).Value = "TTEK"

$ws.Cells.Item(29, 4).Value = 33.19524887145434
$ws.Cells.Item(29, 5).Value = 27.18389511108398
$ws.Cells.Item(29, 6).Value = 34.46275033748218
$ws.Cells.Item(29, 7).Value = 24.73091995175899
$ws.Cells.Item(29, 8).Value = 262773412
$ws.Cells.Item(29, 9).Value = "TTEK"

$ws.Cells.Item(30, 4).Value = 32.23957106398169
$ws.Cells.Item(30, 5).Value = 27.23721694946289
$ws.Cells.Item(30, 6).Value = 33.06482445904211
$ws.Cells.Item(30, 7).Value = 27.10228195818738
$ws.Cells.Item(30, 8).Value = 262773412
$ws.Cells.Item(30, 9).Value = "TTEK"

$ws.Cells.Item(31, 4).Value = 26.85699558085237
$ws.Cells.Item(31, 5).Value = 30.02897262573243
$ws.Cells.Item(31, 6).Value = 30.45216535449418
$ws.Cells.Item(31, 7).Value = 25.75199493972212
$ws.Cells.Item(31, 8).Value = 262773412
$ws.Cells.Item(31, 9).Value = "TTEK"

$ws.Cells.Item(32, 4).Value = 25.34903057818346
$ws.Cells.Item(32, 5).Value = 27.72341728210449
$ws.Cells.Item(32, 6).Value = 27.88432569387681
$ws.Cells.Item(32, 7).Value = 24.1088551425565
$ws.Cells.Item(32, 8).Value = 262773412
$ws.Cells.Item(32, 9).Value = "TTEK"

$ws.Cells.Item(33, 4).Value = 28.8629536415741
$ws.Cells.Item(33, 5).Value = 30.5628547668457
$ws.Cells.Item(33, 6).Value = 30.79474891451654
$ws.Cells.Item(33, 7).Value = 25.8031304915864
$ws.Cells.Item(33, 8).Value = 262773412
$ws.Cells.Item(33, 9).Value = "TTEK"

$ws.Cells.Item(34, 4).Value = 28.81972983750614
$ws.Cells.Item(34, 5).Value = 27.23525428771973
$ws.Cells.Item(34, 6).Value = 29.09725858462779
$ws.Cells.Item(34, 7).Value = 26.25110835659463
$ws.Cells.Item(34, 8).Value = 262773412
$ws.Cells.Item(34, 9).Value = "TTEK"

$ws.Cells.Item(35, 4).Value = 32.0922103700198
$ws.Cells.Item(35, 5).Value = 33.36603927612305
$ws.Cells.Item(35, 6).Value = 34.16661943727939
$ws.Cells.Item(35, 7).Value = 31.89108087575755
$ws.Cells.Item(35, 8).Value = 262773412
$ws.Cells.Item(35, 9).Value = "TTEK"

$ws.Cells.Item(36, 4).Value = 29.86959918886632
$ws.Cells.Item(36, 5).Value = 29.80639457702637
$ws.Cells.Item(36, 6).Value = 33.05545255943224
$ws.Cells.Item(36, 7).Value = 28.31321175742861
$ws.Cells.Item(36, 8).Value = 262773412
$ws.Cells.Item(36, 9).Value = "TTEK"

$ws.Cells.Item(37, 4).Value = 32.73694332215805
$ws.Cells.Item(37, 5).Value = 31.29278182983398
$ws.Cells.Item(37, 6).Value = 33.33241342256853
$ws.Cells.Item(37, 7).Value = 31.26508572159684
$ws.Cells.Item(37, 8).Value = 262773412
$ws.Cells.Item(37, 9).Value = "TTEK"

$ws.Cells.Item(38, 4).Value = 36.57470622286429
$ws.Cells.Item(38, 5).Value = 38.57768630981445
$ws.Cells.Item(38, 6).Value = 38.92637651653513
$ws.Cells.Item(38, 7).Value = 36.16063872796724
$ws.Cells.Item(38, 8).Value = 262773412
$ws.Cells.Item(38, 9).Value = "TTEK"

$ws.Cells.Item(39, 4).Value = 40.66833760930426
$ws.Cells.Item(39, 5).Value = 42.30300521850586
$ws.Cells.Item(39, 6).Value = 42.88823266716404
$ws.Cells.Item(39, 7).Value = 39.10905542317585
$ws.Cells.Item(39, 8).Value = 262773412
$ws.Cells.Item(39, 9).Value = "TTEK"

$ws.Cells.Item(40, 4).Value = 46.76884174008141
$ws.Cells.Item(40, 5).Value = 48.54663467407227
$ws.Cells.Item(40, 6).Value = 50.34428787206557
$ws.Cells.Item(40, 7).Value = 46.11334315626162
$ws.Cells.Item(40, 8).Value = 262773412
$ws.Cells.Item(40, 9).Value = "TTEK"

$ws.Cells.Item(41, 4).Value = 39.8820670207078
$ws.Cells.Item(41, 5).Value = 36.60000228881836
$ws.Cells.Item(41, 6).Value = 42.74641614761383
$ws.Cells.Item(41, 7).Value = 30.92103533518907
$ws.Cells.Item(41, 8).Value = 262773412
$ws.Cells.Item(41, 9).Value = "TTEK"

$ws.Cells.Item(42, 4).Value = 29.14344817806013
$ws.Cells.Item(42, 5).Value = 31.07638168334961
$ws.Cells.Item(42, 6).Value = 31.48488798441175
$ws.Cells.Item(42, 7).Value = 27.1706613728722
$ws.Cells.Item(42, 8).Value = 262773412
$ws.Cells.Item(42, 9).Value = "TTEK"

$ws.Cells.Item(43, 4).Value = 35.88483951065574
$ws.Cells.Item(43, 5).Value = 36.67340850830078
$ws.Cells.Item(43, 6).Value = 38.58992658037779
$ws.Cells.Item(43, 7).Value = 35.57540003102994
$ws.Cells.Item(43, 8).Value = 262773412
$ws.Cells.Item(43, 9).Value = "TTEK"
